# Applies the "feat: add workbook meta and special @TABLEAU sheet" change:
#  - adds a new "Exchange" worksheet (Reward/Exchange bridge table)
#  - adds a new "@TABLEAU" worksheet (sheet/alias/name-line/type-line manifest)
#  - tweaks a couple of pre-existing selections / the A1 style on "Activity"

$wb = $excel.ActiveWorkbook

$activity = $wb.Worksheets.Item("Activity")
$reward   = $wb.Worksheets.Item("Reward")

# ---------------------------------------------------------------------------
# 1. "Activity" sheet: header cell A1 picks up word-wrap (new cellXf, same
#    shared string / value), and the selection moves.
# ---------------------------------------------------------------------------
$activity.Range("A1").WrapText = $true
$activity.Range("H18").Select()

# ---------------------------------------------------------------------------
# 2. "Reward" sheet: just the selection changes.
# ---------------------------------------------------------------------------
$reward.Range("A1:B3").Select()

# ---------------------------------------------------------------------------
# 3. New "Exchange" worksheet, inserted after "Reward".
# ---------------------------------------------------------------------------
$exchange = $wb.Worksheets.Add($null, $reward)
$exchange.Name = "Exchange"

$exchange.Columns.Item(1).ColumnWidth = 19
$exchange.Columns.Item(2).ColumnWidth = 13.75

$exchange.Range("A1").Value = "id" + [char]10 + "ID"
$exchange.Range("B1").Value = "desc" + [char]10 + "Desc"
$exchange.Range("C1").Value = "Date"

$exchange.Range("A2").Value = "INTEGER" + [char]10 + "map<uint32, Reward>"
$exchange.Range("B2").Value = "VARCHAR(64)" + [char]10 + "string"
$exchange.Range("C2").Value = "date"

# Pull in the grey/boxed header cellXf (same one "Activity"/"Reward" row-1
# headers use) before turning word-wrap on, so it lands on the very same
# cellXf the author's workbook uses.
$activity.Range("B1").Copy()
$exchange.Range("A1:C2").PasteSpecial(-4122)
$exchange.Range("A1:C2").WrapText = $true

$exchange.Range("A3").Value = $reward.Range("A3").Text
$exchange.Range("B3").Value = $reward.Range("B3").Text
$exchange.Range("C3").Value = $reward.Range("B3").Text
$reward.Range("A3").Copy()
$exchange.Range("A3:C3").PasteSpecial(-4122)

$exchange.Range("A4").Value = 1
$exchange.Range("B4").Value = "award1"

$exchange.Range("A5").Value = 2
$exchange.Range("B5").Value = "award2"

$exchange.Range("A1").Select()

# ---------------------------------------------------------------------------
# 4. New "@TABLEAU" worksheet, inserted after "Exchange" -- this is the
#    special manifest sheet tableau-style tooling looks for.
# ---------------------------------------------------------------------------
$tableau = $wb.Worksheets.Add($null, $exchange)
$tableau.Name = "@TABLEAU"

$tableau.Columns.Item(1).ColumnWidth = 10.125
$tableau.Columns.Item(2).ColumnWidth = 12.375
$tableau.Columns.Item(3).ColumnWidth = 12.125
$tableau.Columns.Item(4).ColumnWidth = 11.875

$tableau.Range("A1").Value = "Sheet"
$tableau.Range("B1").Value = "Alias"
$tableau.Range("C1").Value = "NameCellLine"
$tableau.Range("D1").Value = "TypeCellLine"

$activity.Range("B1").Copy()
$tableau.Range("A1:D1").PasteSpecial(-4122)
$tableau.Range("A1:B1").WrapText = $true

$tableau.Range("A2").Value = "Activity"
$tableau.Range("B2").Value = "ActivityInfo"
$tableau.Range("C2").Value = $activity.Range("A3").Text
$tableau.Range("D2").Value = $activity.Range("A3").Text

$tableau.Range("A3").Value = "Exchange"
$tableau.Range("B3").Value = "ExchangeInfo"
$tableau.Range("C3").Value = $activity.Range("C3").Text
$tableau.Range("D3").Value = $activity.Range("C3").Text

$tableau.Range("G18").Select()

# ---------------------------------------------------------------------------
# 5. Final bookkeeping: "Exchange" is the tab that ends up active/selected.
# ---------------------------------------------------------------------------
$exchange.Activate()
